$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clarify the PASS/FAIL status formula so it reports "ERROR" whenever the
# actual result (column B) is itself an error, rather than propagating the
# error (or comparing it) and reporting FAIL/#NUM!/#N/A. Filled down over the
# whole status column (D3:D20) using a relative reference, matching how the
# original shared "PASS"/"FAIL" formula was authored.
$ws.Range("D3:D20").Formula = '=IF(ISERROR(B3),"ERROR",IF(ISERROR(C3),"FAIL",IF(B3=C3,"PASS","FAIL")))'
